$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column BV: header "18-sep" plus one data value per existing row (2-18)
$ws.Range("BV1").Value = "18-sep"
$ws.Range("BV2").Value = 0
$ws.Range("BV3").Value = 13.477752678266757
$ws.Range("BV4").Value = 18.228496928146512
$ws.Range("BV5").Value = 13.965362300654338
$ws.Range("BV6").Value = 0
$ws.Range("BV7").Value = 6.4794461216874097
$ws.Range("BV8").Value = 15.171960384450029
$ws.Range("BV9").Value = 7.58655075961574
$ws.Range("BV10").Value = 9.0841000158561442
$ws.Range("BV11").Value = 11.948169857989011
$ws.Range("BV12").Value = 0
$ws.Range("BV13").Value = 10.572743472774256
$ws.Range("BV14").Value = 0
$ws.Range("BV15").Value = 0
$ws.Range("BV16").Value = 15.684717266013774
$ws.Range("BV17").Value = 0
$ws.Range("BV18").Value = 0

# Move the active selection, matching the author's final cursor position
$ws.Range("BX5").Select()
